$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -3.3867
$ws.Range("B4").Value = -0.0086
$ws.Range("B5").Value = 0.0844
$ws.Range("B6").Value = -0.0463
$ws.Range("B7").Value = 0.0588
$ws.Range("B8").Value = -0.1382
$ws.Range("B9").Value = 0.0306
$ws.Range("B10").Value = -0.0874
$ws.Range("B11").Value = 0.1096
$ws.Range("B12").Value = 0.0072
$ws.Range("B13").Value = 0.031
$ws.Range("B14").Value = 0.0035
$ws.Range("B15").Value = -0.194
$ws.Range("B16").Value = 0.0313
$ws.Range("B17").Value = 0.0455
$ws.Range("B18").Value = 0.0006
$ws.Range("B19").Value = -0.05
$ws.Range("B20").Value = 0.0232
$ws.Range("B21").Value = 0.8692
$ws.Range("B22").Value = 0.1137
$ws.Range("B23").Value = 0.0003
$ws.Range("B24").Value = 0.1188
$ws.Range("B25").Value = 0.12
